# Update "minhkhoi" row (row 4) match-history stats after another match played.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 600
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = ";0;0;0;1;1;0"
$ws.Range("G4").Value = ";0;0;0;0;3;42"
$ws.Range("H4").Value = ";-100;-200;-100;+600;+700;-200"
